# Refresh the cryptos list with the latest GitHub Actions scrape:
#  - Stellar and ONDO swapped ranking positions (rows 49-50), bringing
#    each coin's own price / 1h-volume-change figures along with it.
#  - Every other row keeps its coin/link but gets updated Price (D) and
#    Volume(1h) (E) figures.
#
# Column D holds plain text (several prices use a "thousands." style like
# "67.921.71" that Excel can't parse as a number anyway, but some, like
# "582.12", look numeric). To keep every D cell a literal text value
# (matching the original inlineStr cells) instead of letting Excel's
# auto-type-detection turn the numeric-looking ones into real numbers,
# the whole D2:D51 range is temporarily switched to Text format while the
# values are written, then restored to the default "Normal" style so the
# saved file carries no stray per-cell formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Rows 49 and 50 swapped coin identity (Stellar <-> ONDO) along with
# their price / volume figures.
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +5.49%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.103"
$ws.Range("E50").Value = "  +1.27%  "

# Refreshed price / 1h volume-change figures for the rest of the sheet.
$ws.Range("D2").Value = "67.921.71"
$ws.Range("E2").Value = "  +3.32%  "

$ws.Range("D3").Value = "3.283.18"
$ws.Range("E3").Value = "  +3.51%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "582.12"
$ws.Range("E5").Value = "  +1.85%  "

$ws.Range("D6").Value = "183.46"
$ws.Range("E6").Value = "  +6.75%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  +0.86%  "

$ws.Range("E9").Value = "  +7.71%  "

$ws.Range("D10").Value = "6.73"
$ws.Range("E10").Value = "  +1.75%  "

$ws.Range("E11").Value = "  +6.25%  "

$ws.Range("D12").Value = "3.852.18"
$ws.Range("E12").Value = "  +3.44%  "

$ws.Range("E13").Value = "  +1.47%  "

$ws.Range("D14").Value = "28.79"
$ws.Range("E14").Value = "  +5.78%  "

$ws.Range("D15").Value = "67.899.47"
$ws.Range("E15").Value = "  +3.36%  "

$ws.Range("E16").Value = "  +3.79%  "

$ws.Range("D17").Value = "3.283.57"
$ws.Range("E17").Value = "  +3.48%  "

$ws.Range("E18").Value = "  +2.10%  "

$ws.Range("D19").Value = "13.57"
$ws.Range("E19").Value = "  +5.13%  "

$ws.Range("D20").Value = "377.26"
$ws.Range("E20").Value = "  +4.21%  "

$ws.Range("E21").Value = "  +5.82%  "

$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("D23").Value = "71.35"
$ws.Range("E23").Value = "  +3.50%  "

$ws.Range("D24").Value = "0.515"
$ws.Range("E24").Value = "  +3.91%  "

$ws.Range("E25").Value = "  +6.01%  "

$ws.Range("D26").Value = "9.81"
$ws.Range("E26").Value = "  -0.40%  "

$ws.Range("E27").Value = "  +2.77%  "

$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("E29").Value = "  +3.25%  "

$ws.Range("D30").Value = "5.73"
$ws.Range("E30").Value = "  +6.35%  "

$ws.Range("D31").Value = "22.95"
$ws.Range("E31").Value = "  +4.13%  "

$ws.Range("E32").Value = "  +7.35%  "

$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("D34").Value = "6.98"
$ws.Range("E34").Value = "  +5.27%  "

$ws.Range("E35").Value = "  +5.65%  "

$ws.Range("D36").Value = "162.00"
$ws.Range("E36").Value = "  +0.45%  "

$ws.Range("E37").Value = "  +2.20%  "

$ws.Range("E38").Value = "  +2.59%  "

$ws.Range("D39").Value = "27.06"
$ws.Range("E39").Value = "  +2.68%  "

$ws.Range("D40").Value = "6.79"
$ws.Range("E40").Value = "  +10.01%  "

$ws.Range("E41").Value = "  +10.49%  "

$ws.Range("E42").Value = "  +5.38%  "

$ws.Range("D43").Value = "25.93"
$ws.Range("E43").Value = "  +9.17%  "

$ws.Range("D44").Value = "352.60"
$ws.Range("E44").Value = "  +7.22%  "

$ws.Range("D45").Value = "2.667.78"
$ws.Range("E45").Value = "  +0.72%  "

$ws.Range("E46").Value = "  +2.96%  "

$ws.Range("E47").Value = "  +3.52%  "

$ws.Range("E48").Value = "  +4.13%  "

$ws.Range("D51").Value = "31.16"
$ws.Range("E51").Value = "  +3.45%  "

# Restore the default style on the price column so no residual
# number-format style is left behind on any cell.
$priceRange.Style = "Normal"
